# Holly added "S.GISH" as the harvester for the bioSamples / rnaSamples
# sheet. Update the "harvester" column (column B) for every data row
# (rows 2-25) from "Retrofitted_0759" to "S.GISH".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2:B25").Value = "S.GISH"
